# Auto-generated edit script applying numeric updates to Kraken_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4 (item id 5470) on sheet ALC
$ws.Range("H4").Value = 3083.2856
$ws.Range("I4").Value = 3588.8333
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 3588.8333
$ws.Range("L4").Value = 50
$ws.Range("M4").Value = -3474.8333
$ws.Range("N4").Value = -278

# Row 6 (item id 4564) on sheet ALC
$ws.Range("H6").Value = 33.857143
$ws.Range("I6").Value = 7.4
$ws.Range("K6").Value = 22.2
$ws.Range("M6").Value = 89.8

# Row 19 (item id 7015) on sheet ALC
$ws.Range("H19").Value = 620.55554
$ws.Range("I19").Value = 698.25
$ws.Range("J19").Value = 558.4
$ws.Range("K19").Value = 698.25
$ws.Range("L19").Value = 558.4
$ws.Range("M19").Value = -523.25
$ws.Range("N19").Value = -908.4

# Row 33 (item id 5512) on sheet ALC
$ws.Range("H33").Value = 85.5625
$ws.Range("I33").Value = 91
$ws.Range("J33").Value = 47.5
$ws.Range("K33").Value = 91
$ws.Range("L33").Value = 47.5
$ws.Range("M33").Value = 138
$ws.Range("N33").Value = -505.5

# Row 69 (item id 12616) on sheet ALC
$ws.Range("H69").Value = 500
$ws.Range("I69").Value = 500
$ws.Range("K69").Value = 1500
$ws.Range("M69").Value = -626

# Row 72 (item id 12616) on sheet ALC
$ws.Range("H72").Value = 500
$ws.Range("I72").Value = 500
$ws.Range("K72").Value = 4500
$ws.Range("M72").Value = -132

# Row 92 (item id 19901) on sheet ALC
$ws.Range("H92").Value = 1872.7142
$ws.Range("I92").Value = 1851.6666
$ws.Range("K92").Value = 1851.6666
$ws.Range("M92").Value = -603.6666

# Row 116 (item id 27778) on sheet ALC
$ws.Range("H116").Value = 3450
$ws.Range("I116").Value = 3500
$ws.Range("J116").Value = 3400
$ws.Range("K116").Value = 3500
$ws.Range("L116").Value = 3400
$ws.Range("M116").Value = -58
$ws.Range("N116").Value = -10284

# Row 137 (item id 44013) on sheet ALC
$ws.Range("H137").Value = 3482.647
$ws.Range("J137").Value = 3726.923
$ws.Range("L137").Value = 11180.769
$ws.Range("N137").Value = -16280.769

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (item id 27713) on sheet ARM
$ws.Range("H2").Value = 603.6667
$ws.Range("I2").Value = 603.6667
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 603.6667
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -490.6667
$ws.Range("N2").ClearContents()

# Row 32 (item id 44147) on sheet ARM
$ws.Range("H32").Value = 3343.8572
$ws.Range("I32").Value = 2966.4614
$ws.Range("J32").Value = 8250
$ws.Range("K32").Value = 2966.4614
$ws.Range("L32").Value = 8250
$ws.Range("M32").Value = -2679.4614
$ws.Range("N32").Value = -8824

# Row 45 (item id 27714) on sheet ARM
$ws.Range("H45").Value = 2577.92
$ws.Range("I45").Value = 2415.348
$ws.Range("J45").Value = 4447.5
$ws.Range("K45").Value = 2415.348
$ws.Range("L45").Value = 4447.5
$ws.Range("M45").Value = -2038.348
$ws.Range("N45").Value = -5201.5

# Row 61 (item id 43999) on sheet ARM
$ws.Range("H61").Value = 2999
$ws.Range("I61").Value = 2999
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2999
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2787
$ws.Range("N61").ClearContents()

# Row 63 (item id 12528) on sheet ARM
$ws.Range("H63").Value = 1575.5555
$ws.Range("I63").Value = 1547.1666
$ws.Range("J63").Value = 1632.3334
$ws.Range("K63").Value = 1547.1666
$ws.Range("L63").Value = 1632.3334
$ws.Range("M63").Value = -861.1666
$ws.Range("N63").Value = -3004.3334

# Row 66 (item id 12528) on sheet ARM
$ws.Range("H66").Value = 1575.5555
$ws.Range("I66").Value = 1547.1666
$ws.Range("J66").Value = 1632.3334
$ws.Range("K66").Value = 7735.833000000001
$ws.Range("L66").Value = 8161.666999999999
$ws.Range("M66").Value = -4303.833000000001
$ws.Range("N66").Value = -15025.667

# Row 97 (item id 19941) on sheet ARM
$ws.Range("H97").Value = 1894.6
$ws.Range("I97").Value = 849.5714
$ws.Range("J97").Value = 4333
$ws.Range("K97").Value = 849.5714
$ws.Range("L97").Value = 4333
$ws.Range("M97").Value = -353.5714
$ws.Range("N97").Value = -5325

# Row 116 (item id 27713) on sheet ARM
$ws.Range("H116").Value = 603.6667
$ws.Range("I116").Value = 603.6667
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 603.6667
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1690.3333
$ws.Range("N116").ClearContents()

# Row 136 (item id 43999) on sheet ARM
$ws.Range("H136").Value = 2999
$ws.Range("I136").Value = 2999
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8997
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6447
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (item id 27713) on sheet BSM
$ws.Range("H3").Value = 603.6667
$ws.Range("I3").Value = 603.6667
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 603.6667
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -489.6667
$ws.Range("N3").ClearContents()

# Row 105 (item id 19947) on sheet BSM
$ws.Range("H105").Value = 3964.2
$ws.Range("I105").Value = 3964.2
$ws.Range("K105").Value = 3964.2
$ws.Range("M105").Value = -2217.2

# Row 112 (item id 25788) on sheet BSM
$ws.Range("H112").Value = 20434.5
$ws.Range("J112").Value = 20434.5
$ws.Range("L112").Value = 20434.5
$ws.Range("N112").Value = -23388.5

# Row 134 (item id 43998) on sheet BSM
$ws.Range("H134").Value = 12000
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (item id 5361) on sheet CRP
$ws.Range("H7").Value = 168.55556
$ws.Range("I7").Value = 89.40000000000001
$ws.Range("J7").Value = 267.5
$ws.Range("K7").Value = 89.40000000000001
$ws.Range("L7").Value = 267.5
$ws.Range("M7").Value = 23.59999999999999
$ws.Range("N7").Value = -493.5

# Row 22 (item id 5367) on sheet CRP
$ws.Range("H22").Value = 224
$ws.Range("I22").Value = 224
$ws.Range("K22").Value = 224
$ws.Range("M22").Value = 126

# Row 58 (item id 44021) on sheet CRP
$ws.Range("H58").Value = 2900
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 2900
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 2900
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -3306

# Row 125 (item id 34297) on sheet CRP
$ws.Range("H125").Value = 45000
$ws.Range("J125").Value = 45000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -49920

# Row 136 (item id 44021) on sheet CRP
$ws.Range("H136").Value = 2900
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 2900
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 8700
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -13800

$ws = $wb.Worksheets.Item("CUL")
# Row 110 (item id 27857) on sheet CUL
$ws.Range("H110").Value = 1000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 3000
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -11180

# Row 137 (item id 44088) on sheet CUL
$ws.Range("H137").Value = 2000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 6000
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -16200

$ws = $wb.Worksheets.Item("GSM")
# Row 9 (item id 1683) on sheet GSM
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

# Row 70 (item id 14146) on sheet GSM
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# Row 73 (item id 14146) on sheet GSM
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# Row 107 (item id 27802) on sheet GSM
$ws.Range("H107").Value = 610
$ws.Range("I107").Value = 218.2
$ws.Range("J107").Value = 1263
$ws.Range("K107").Value = 218.2
$ws.Range("L107").Value = 1263
$ws.Range("M107").Value = 1701.8
$ws.Range("N107").Value = -5103

# Row 132 (item id 44008) on sheet GSM
$ws.Range("H132").Value = 5172.4
$ws.Range("I132").Value = 5172.4
$ws.Range("K132").Value = 15517.2
$ws.Range("M132").Value = -12987.2

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (item id 36249) on sheet LTW
$ws.Range("H7").Value = 5020.5713
$ws.Range("I7").Value = 4629
$ws.Range("K7").Value = 4629
$ws.Range("M7").Value = -4517

# Row 122 (item id 36247) on sheet LTW
$ws.Range("H122").Value = 5171.75
$ws.Range("I122").Value = 6344.5
$ws.Range("K122").Value = 19033.5
$ws.Range("M122").Value = -16583.5

# Row 126 (item id 36249) on sheet LTW
$ws.Range("H126").Value = 5020.5713
$ws.Range("I126").Value = 4629
$ws.Range("K126").Value = 13887
$ws.Range("M126").Value = -11417

# Row 136 (item id 44060) on sheet LTW
$ws.Range("H136").Value = 3657.6667
$ws.Range("I136").Value = 3657.6667
$ws.Range("K136").Value = 10973.0001
$ws.Range("M136").Value = -8423.000100000001

$ws = $wb.Worksheets.Item("WVR")
# Row 136 (item id 44031) on sheet WVR
$ws.Range("H136").Value = 2040.4286
$ws.Range("I136").Value = 1850.8462
$ws.Range("K136").Value = 5552.5386
$ws.Range("M136").Value = -3002.5386

Write-Host "Applied all Kraken_Profits updates"